$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.776.41"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.87%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.408.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.60%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.411.93"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.96%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.545"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.35"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.122"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.431"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.005.81"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.18%  "
$ws.Range("E14").Value = "  -3.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000192"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.86%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.788.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.80%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.443.58"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "380.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.995"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "71.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.529"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000119"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +24.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.09"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.18%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.41"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.77"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.90%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.984.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0754"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0312"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.98%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.74"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.761"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.43%  "
$ws.Range("E48").Value = "  +3.10%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.18"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +21.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.830"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.12%  "
